# Update Document Index sheet: change Document No. (column A) values
# from VP-NCC-R-001-00x to VP-NCC-R-004-00x for rows 1-5.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "VP-NCC-R-004-001"
$ws.Range("A2").Value = "VP-NCC-R-004-002"
$ws.Range("A3").Value = "VP-NCC-R-004-003"
$ws.Range("A4").Value = "VP-NCC-R-004-004"
$ws.Range("A5").Value = "VP-NCC-R-004-005"
